$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '70.151.69', '  -0.75%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.576.85', '  -1.18%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.01%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '576.14', '  -1.78%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '186.30', '  -3.57%  '),
    @(7, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.569.54', '  -1.26%  '),
    @(8, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.621', '  -2.70%  '),
    @(9, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.08%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.185', '  +1.81%  '),
    @(11, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.652', '  -2.35%  '),
    @(12, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '54.23', '  -5.43%  '),
    @(13, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000302', '  -1.15%  '),
    @(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '9.58', '  -2.55%  '),
    @(15, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '4.160.76', '  -0.96%  '),
    @(16, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '19.54', '  -3.91%  '),
    @(17, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.581.83', '  -1.12%  '),
    @(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '70.157.37', '  -0.74%  '),
    @(19, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '12.48', '  -1.43%  '),
    @(20, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.120', '  -1.17%  '),
    @(21, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.05', '  -0.43%  '),
    @(22, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '500.78', '  +3.10%  '),
    @(23, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '19.54', '  +1.15%  '),
    @(24, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '4.91', '  -3.66%  '),
    @(25, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '4.36', '  -1.90%  '),
    @(26, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '95.05', '  +5.36%  '),
    @(27, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '11.46', '  +0.99%  '),
    @(28, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.96', '  -4.97%  '),
    @(29, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '9.33', '  -1.31%  '),
    @(30, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '31.73', '  -2.44%  '),
    @(31, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '7.62', '  -3.07%  '),
    @(32, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '12.54', '  +2.40%  '),
    @(33, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '66.02', '  -0.89%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.116', '  -8.95%  '),
    @(35, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '571.32', '  -6.24%  '),
    @(36, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '3.17', '  +6.15%  '),
    @(37, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '38.64', '  -3.60%  '),
    @(38, 'TheGraph', 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt', '0.408', '  -0.01%  '),
    @(39, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.999', '  -0.03%  '),
    @(40, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0788', '  -5.11%  '),
    @(41, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '3.31', '  +3.57%  '),
    @(42, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '3.43', '  -3.41%  '),
    @(43, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.136', '  -7.72%  '),
    @(44, 'ThetaToken', 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta', '3.02', '  -3.42%  '),
    @(45, 'ApeXProtocol', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex', '3.56', '  +5.65%  '),
    @(46, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0449', '  -1.01%  '),
    @(47, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '3.205.72', '  -3.33%  '),
    @(48, 'THORChain', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune', '9.40', '  -3.05%  '),
    @(49, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.136', '  -1.76%  '),
    @(50, 'OceanProtocol', 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean', '1.49', '  +24.01%  '),
    @(51, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  +0.00%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]
    $eCell.Style = "Normal"
}
